# Trade #39 closed at 2026-02-17 12:41:53 - unknown UNKNOWN +0.000%
$wb = $excel.ActiveWorkbook

# --- Summary sheet ---
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B3").Value = 1200.3
$wsSummary.Range("B4").Value = 0.29
$wsSummary.Range("B5").Value = 0.15
$wsSummary.Range("B6").Value = 39
$wsSummary.Range("B8").Value = 15
$wsSummary.Range("B9").Value = 38.46

# --- Strategy Status sheet (MarketMaking row) ---
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("C4").Value = 100.3
$wsStatus.Range("D4").Value = 39
$wsStatus.Range("E4").Value = 0.29
$wsStatus.Range("F4").Value = 0.3
$wsStatus.Range("G4").Value = 38.46

# --- All Trades & MarketMaking sheets: update trade #39 (row 40) to CLOSED ---
$tradeSheets = @("All Trades", "MarketMaking")
foreach ($name in $tradeSheets) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("G40").Value = 0.53
    $ws.Range("H40").Value = "CLOSED"
    $ws.Range("I40").Value = -45.3608
    $ws.Range("J40").Value = -0.44
    $ws.Range("K40").Value = 100.3
    $ws.Range("P40").Value = "early_exit"
    $ws.Range("Q40").Value = 2.57
}
